# Updated symbol list on Thu Dec 15 22:38:32 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" column (D) for the crypto rows whose quotes moved.
# Values are written with a leading apostrophe so Excel keeps storing them
# as literal text (matching the original inline-string cells and
# preserving exact formatting such as trailing zeros) instead of
# re-interpreting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula  = "'257.69"
$ws.Range("D3").Formula  = "'22.80"
$ws.Range("D4").Formula  = "'6.164"
$ws.Range("D5").Formula  = "'0.06067"
$ws.Range("D6").Formula  = "'6.718"
$ws.Range("D7").Formula  = "'3.454"
$ws.Range("D8").Formula  = "'1.357"
$ws.Range("D9").Formula  = "'0.7966"
$ws.Range("D10").Formula = "'0.1583"
$ws.Range("D11").Formula = "'0.08041"
$ws.Range("D12").Formula = "'0.03353"
$ws.Range("D13").Formula = "'0.03084"
$ws.Range("D14").Formula = "'0.09301"
$ws.Range("D15").Formula = "'3.895"
$ws.Range("D16").Formula = "'0.001693"
$ws.Range("D17").Formula = "'0.04840"
$ws.Range("D18").Formula = "'0.0006157"
$ws.Range("D19").Formula = "'0.006215"
$ws.Range("D20").Formula = "'0.001102"
$ws.Range("D21").Formula = "'0.003378"
$ws.Range("D22").Formula = "'0.0001502"
$ws.Range("D23").Formula = "'3.686"
$ws.Range("D24").Formula = "'2.261"
$ws.Range("D26").Formula = "'0.1226"
$ws.Range("D27").Formula = "'0.0003019"
$ws.Range("D40").Formula = "'0.04574"
$ws.Range("D41").Formula = "'0.007140"
$ws.Range("D42").Formula = "'0.003905"
$ws.Range("D44").Formula = "'0.009931"
$ws.Range("D45").Formula = "'0.002974"
$ws.Range("D46").Formula = "'0.00005975"
$ws.Range("D48").Formula = "'0.7509"
$ws.Range("D49").Formula = "'0.1072"
$ws.Range("D50").Formula = "'0.00001502"
$ws.Range("D51").Formula = "'0.01011"
